# Error Calculations and Plots
# Applies the edits to the "missing data" table:
#  - removes the "RM 232" row (old row 26) and the "SC 92" row (old row 28,
#    which becomes row 27 once the first row is removed), shifting every
#    subsequent row up so the table now spans A1:F33 instead of A1:F35
#  - updates a handful of individual cell values (some numbers become blank /
#    missing, and some previously-blank cells now contain numbers) to match
#    the newly re-sampled "missing data" selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that disappear from the table entirely.
# Deleting row 26 first shifts the old "SC 92" row up to row 27, so we
# delete row 27 next to remove it as well.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Individual cell edits on the resulting (renumbered) sheet.
$ws.Range("F5").ClearContents()
$ws.Range("E6").Value = -5.7
$ws.Range("E8").ClearContents()
$ws.Range("F11").Value = 17.65
$ws.Range("E19").Value = -6.5
$ws.Range("F19").ClearContents()
$ws.Range("E21").ClearContents()
$ws.Range("E23").Value = -7
$ws.Range("F23").Value = 16.48
$ws.Range("F25").Value = 16.6
$ws.Range("D26").ClearContents()
$ws.Range("D27").Value = -14.6
$ws.Range("E27").ClearContents()
$ws.Range("F27").ClearContents()
$ws.Range("D29").ClearContents()
$ws.Range("E29").Value = -6.8
$ws.Range("F29").ClearContents()
$ws.Range("F30").Value = 16.89
$ws.Range("F33").Value = 17.53
